$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Update title / date shared-string runs (surgical, in place) ----
# "Volume 32   Number  18" -> "...19"  (cell A8)
$ws.Range("A8").Characters(21,2).Text = "19"

# "Report Covering the Week  4/28/2025  Through  5/4/2025" -> "...5/5/2025...5/11/2025" (cell C9)
# Edit right-to-left so earlier character offsets stay valid.
$ws.Range("C9").Characters(47,8).Text = "5/11/2025"
$ws.Range("C9").Characters(27,9).Text = "5/5/2025"

# ---- Helper: convert a numeric cell to shared-text "n/a" style (s=13), pulling the
# exact style from a donor cell that already carries it, so no new style gets created ----
function Set-NACell($addr, $text, $styleDonor) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
    $ws.Range($styleDonor).Copy()
    $ws.Range($addr).PasteSpecial(-4122)
}

# Helper: convert a shared-text "n/a" cell back into a real number, reusing the pre-existing
# number-format style (this matches the style already in the stylesheet, so Excel reuses it).
function Set-NumCell($addr, $value, $numFmt) {
    $ws.Range($addr).NumberFormat = $numFmt
    $ws.Range($addr).Value = $value
}

# ---- Cell value / type updates ----
$ws.Range("N15").Value = -55.555555555555
Set-NACell "C16" "0" "C14"
$ws.Range("E16").Value = -100
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = -37.5
$ws.Range("I16").Value = 27
$ws.Range("J16").Value = 33
$ws.Range("K16").Value = -18.181818181818
$ws.Range("L16").Value = -22.857142857142
$ws.Range("M16").Value = -15.625
$ws.Range("N16").Value = -87.557603686635
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = -87.5
$ws.Range("F17").Value = 6
$ws.Range("G17").Value = 14
$ws.Range("H17").Value = -57.142857142857
$ws.Range("I17").Value = 25
$ws.Range("J17").Value = 43
$ws.Range("K17").Value = -41.860465116279
$ws.Range("L17").Value = -44.444444444444
$ws.Range("M17").Value = -7.407407407407
$ws.Range("N17").Value = -28.571428571428
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -80
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = -45.454545454545
$ws.Range("I18").Value = 53
$ws.Range("J18").Value = 62
$ws.Range("K18").Value = -14.516129032258
$ws.Range("L18").Value = 3.921568627450
$ws.Range("M18").Value = 26.190476190476
$ws.Range("N18").Value = -81.138790035587
$ws.Range("C19").Value = 16
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = 100
$ws.Range("F19").Value = 59
$ws.Range("G19").Value = 54
$ws.Range("H19").Value = 9.259259259259
$ws.Range("I19").Value = 242
$ws.Range("J19").Value = 253
$ws.Range("K19").Value = -4.347826086956
$ws.Range("L19").Value = -3.968253968253
$ws.Range("M19").Value = -1.224489795918
$ws.Range("N19").Value = -65.625
Set-NACell "C20" "0" "C14"
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 5
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = -16.666666666666
$ws.Range("J20").Value = 15
$ws.Range("K20").Value = -13.333333333333
$ws.Range("L20").Value = -56.666666666666
$ws.Range("M20").Value = 30
$ws.Range("N20").Value = -96.969696969697
$ws.Range("C21").Value = 18
$ws.Range("D21").Value = 23
$ws.Range("E21").Value = -21.739130434782
$ws.Range("F21").Value = 82
$ws.Range("G21").Value = 93
$ws.Range("H21").Value = -11.827956989247
$ws.Range("I21").Value = 364
$ws.Range("J21").Value = 408
$ws.Range("K21").Value = -10.784313725490
$ws.Range("L21").Value = -13.126491646778
$ws.Range("M21").Value = 1.111111111111
$ws.Range("N21").Value = -78.281622911694
$ws.Range("D23").Value = 5
$ws.Range("F23").Value = 1
$ws.Range("G23").Value = 7
$ws.Range("H23").Value = -85.714285714285
$ws.Range("J23").Value = 20
$ws.Range("K23").Value = -40
$ws.Range("L23").Value = -33.333333333333
$ws.Range("D24").Value = 30
$ws.Range("E24").Value = 13.333333333333
$ws.Range("F24").Value = 130
$ws.Range("G24").Value = 127
$ws.Range("H24").Value = 2.362204724409
$ws.Range("I24").Value = 469
$ws.Range("J24").Value = 431
$ws.Range("K24").Value = 8.816705336426
$ws.Range("L24").Value = 23.097112860892
$ws.Range("M24").Value = 35.549132947976
$ws.Range("C25").Value = 24
$ws.Range("E25").Value = 20
$ws.Range("F25").Value = 80
$ws.Range("G25").Value = 92
$ws.Range("H25").Value = -13.043478260869
$ws.Range("I25").Value = 360
$ws.Range("J25").Value = 335
$ws.Range("K25").Value = 7.462686567164
$ws.Range("L25").Value = 18.811881188118
$ws.Range("C26").Value = 3
$ws.Range("D26").Value = 11
$ws.Range("E26").Value = -72.727272727272
$ws.Range("F26").Value = 13
$ws.Range("G26").Value = 25
$ws.Range("H26").Value = -48
$ws.Range("I26").Value = 76
$ws.Range("J26").Value = 81
$ws.Range("K26").Value = -6.172839506172
$ws.Range("L26").Value = -5
$ws.Range("M26").Value = -17.391304347826
Set-NACell "C27" "0" "C14"
$ws.Range("L27").Value = -12.5
$ws.Range("C28").Value = 1
Set-NACell "D28" "0" "C14"
Set-NACell "E28" "***.*" "E14"
Set-NumCell "D29" 1 "#,##0"
Set-NumCell "E29" -100 "#,##0.0;""-""#,##0.0"
Set-NumCell "G29" 1 "#,##0"
Set-NumCell "H29" -100 "#,##0.0;""-""#,##0.0"
Set-NumCell "J29" 1 "#,##0"
Set-NumCell "K29" -100 "#,##0.0;""-""#,##0.0"
Set-NumCell "D30" 1 "#,##0"
Set-NumCell "E30" -100 "#,##0.0;""-""#,##0.0"
Set-NumCell "G30" 1 "#,##0"
Set-NumCell "H30" -100 "#,##0.0;""-""#,##0.0"
Set-NumCell "J30" 1 "#,##0"
Set-NumCell "K30" -100 "#,##0.0;""-""#,##0.0"
Set-NACell "D31" "0" "C14"
Set-NACell "E31" "***.*" "E14"
Set-NumCell "L31" 200 "#,##0.0;""-""#,##0.0"
